# Update NATMI ligand-receptor TPM-derived statistics (Cd34-Sell) with new TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 260.1928506666666
$ws.Range("H2").Value = 780.5785519999999
$ws.Range("I2").Value = 0.8191301249666345
$ws.Range("J2").Value = 0.8191301249666346
$ws.Range("M2").Value = 0.353079
$ws.Range("N2").Value = 1.059237
$ws.Range("O2").Value = 0.01390801122570493
$ws.Range("P2").Value = 0.01390801122570493
$ws.Range("Q2").Value = 91.86863152053598
$ws.Range("R2").Value = 826.8176836848239
$ws.Range("S2").Value = 0.01139247097334903
$ws.Range("T2").Value = 0.01139247097334904

# Row 3
$ws.Range("G3").Value = 260.1928506666666
$ws.Range("H3").Value = 780.5785519999999
$ws.Range("I3").Value = 0.8191301249666345
$ws.Range("J3").Value = 0.8191301249666346
$ws.Range("M3").Value = 24.359095
$ws.Range("N3").Value = 73.077285
$ws.Range("O3").Value = 0.9595205795530543
$ws.Range("P3").Value = 0.9595205795530543
$ws.Range("Q3").Value = 6338.062367710146
$ws.Range("R3").Value = 57042.56130939132
$ws.Range("S3").Value = 0.7859722122373509
$ws.Range("T3").Value = 0.785972212237351

# Row 4
$ws.Range("G4").Value = 260.1928506666666
$ws.Range("H4").Value = 780.5785519999999
$ws.Range("I4").Value = 0.8191301249666345
$ws.Range("J4").Value = 0.8191301249666346
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6745613333333332
$ws.Range("N4").Value = 2.023684
$ws.Range("O4").Value = 0.02657140922124081
$ws.Range("P4").Value = 0.02657140922124081
$ws.Range("Q4").Value = 175.5160362695075
$ws.Range("R4").Value = 1579.644326425568
$ws.Range("S4").Value = 0.02176544175593457
$ws.Range("T4").Value = 0.02176544175593457

# Row 5
$ws.Range("I5").Value = 0.1763970508574364
$ws.Range("J5").Value = 0.1763970508574364
$ws.Range("M5").Value = 0.353079
$ws.Range("N5").Value = 1.059237
$ws.Range("O5").Value = 0.01390801122570493
$ws.Range("P5").Value = 0.01390801122570493
$ws.Range("Q5").Value = 19.783615780449
$ws.Range("R5").Value = 178.052542024041
$ws.Range("S5").Value = 0.002453332163506469
$ws.Range("T5").Value = 0.002453332163506469

# Row 6
$ws.Range("I6").Value = 0.1763970508574364
$ws.Range("J6").Value = 0.1763970508574364
$ws.Range("M6").Value = 24.359095
$ws.Range("N6").Value = 73.077285
$ws.Range("O6").Value = 0.9595205795530543
$ws.Range("P6").Value = 0.9595205795530543
$ws.Range("Q6").Value = 1364.881446473612
$ws.Range("R6").Value = 12283.93301826251
$ws.Range("S6").Value = 0.169256600470177
$ws.Range("T6").Value = 0.169256600470177

# Row 7
$ws.Range("I7").Value = 0.1763970508574364
$ws.Range("J7").Value = 0.1763970508574364
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6745613333333332
$ws.Range("N7").Value = 2.023684
$ws.Range("O7").Value = 0.02657140922124081
$ws.Range("P7").Value = 0.02657140922124081
$ws.Range("Q7").Value = 37.79681668695689
$ws.Range("R7").Value = 340.171350182612
$ws.Range("S7").Value = 0.004687118223752971
$ws.Range("T7").Value = 0.004687118223752971

# Row 8
$ws.Range("G8").Value = 0.6130636666666667
$ws.Range("H8").Value = 1.839191
$ws.Range("I8").Value = 0.00193002581201784
$ws.Range("J8").Value = 0.00193002581201784
$ws.Range("M8").Value = 0.353079
$ws.Range("N8").Value = 1.059237
$ws.Range("O8").Value = 0.01390801122570493
$ws.Range("P8").Value = 0.01390801122570493
$ws.Range("Q8").Value = 0.216459906363
$ws.Range("R8").Value = 1.948139157267
$ws.Range("S8").Value = 0.00002684282065944439
$ws.Range("T8").Value = 0.00002684282065944439

# Row 9
$ws.Range("G9").Value = 0.6130636666666667
$ws.Range("H9").Value = 1.839191
$ws.Range("I9").Value = 0.00193002581201784
$ws.Range("J9").Value = 0.00193002581201784
$ws.Range("M9").Value = 24.359095
$ws.Range("N9").Value = 73.077285
$ws.Range("O9").Value = 0.9595205795530543
$ws.Range("P9").Value = 0.9595205795530543
$ws.Range("Q9").Value = 14.93367609738167
$ws.Range("R9").Value = 134.403084876435
$ws.Range("S9").Value = 0.001851899485699712
$ws.Range("T9").Value = 0.001851899485699712

# Row 10
$ws.Range("G10").Value = 0.6130636666666667
$ws.Range("H10").Value = 1.839191
$ws.Range("I10").Value = 0.00193002581201784
$ws.Range("J10").Value = 0.00193002581201784
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6745613333333332
$ws.Range("N10").Value = 2.023684
$ws.Range("O10").Value = 0.02657140922124081
$ws.Range("P10").Value = 0.02657140922124081
$ws.Range("Q10").Value = 0.4135490444048888
$ws.Range("R10").Value = 3.721941399644
$ws.Range("S10").Value = 0.00005128350565868362
$ws.Range("T10").Value = 0.00005128350565868362

# Row 11
$ws.Range("G11").Value = 0.5097843333333333
$ws.Range("H11").Value = 1.529353
$ws.Range("I11").Value = 0.001604885390199778
$ws.Range("J11").Value = 0.001604885390199778
$ws.Range("M11").Value = 0.353079
$ws.Range("N11").Value = 1.059237
$ws.Range("O11").Value = 0.01390801122570493
$ws.Range("P11").Value = 0.01390801122570493
$ws.Range("Q11").Value = 0.179994142629
$ws.Range("R11").Value = 1.619947283661
$ws.Range("S11").Value = 0.00002232076402286834
$ws.Range("T11").Value = 0.00002232076402286835

# Row 12
$ws.Range("G12").Value = 0.5097843333333333
$ws.Range("H12").Value = 1.529353
$ws.Range("I12").Value = 0.001604885390199778
$ws.Range("J12").Value = 0.001604885390199778
$ws.Range("M12").Value = 24.359095
$ws.Range("N12").Value = 73.077285
$ws.Range("O12").Value = 0.9595205795530543
$ws.Range("P12").Value = 0.9595205795530543
$ws.Range("Q12").Value = 12.41788500517833
$ws.Range("R12").Value = 111.760965046605
$ws.Range("S12").Value = 0.00153992055972072
$ws.Range("T12").Value = 0.00153992055972072

# Row 13
$ws.Range("G13").Value = 0.5097843333333333
$ws.Range("H13").Value = 1.529353
$ws.Range("I13").Value = 0.001604885390199778
$ws.Range("J13").Value = 0.001604885390199778
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6745613333333332
$ws.Range("N13").Value = 2.023684
$ws.Range("O13").Value = 0.02657140922124081
$ws.Range("P13").Value = 0.02657140922124081
$ws.Range("Q13").Value = 0.3438807996057777
$ws.Range("R13").Value = 3.094927196452
$ws.Range("S13").Value = 0.00004264406645618904
$ws.Range("T13").Value = 0.00004264406645618904

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2979236666666666
$ws.Range("H14").Value = 0.893771
$ws.Range("I14").Value = 0.0009379129737112659
$ws.Range("J14").Value = 0.000937912973711266
$ws.Range("M14").Value = 0.353079
$ws.Range("N14").Value = 1.059237
$ws.Range("O14").Value = 0.01390801122570493
$ws.Range("P14").Value = 0.01390801122570493
$ws.Range("Q14").Value = 0.105190590303
$ws.Range("R14").Value = 0.9467153127269999
$ws.Range("S14").Value = 0.00001304450416711058
$ws.Range("T14").Value = 0.00001304450416711058

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2979236666666666
$ws.Range("H15").Value = 0.893771
$ws.Range("I15").Value = 0.0009379129737112659
$ws.Range("J15").Value = 0.000937912973711266
$ws.Range("M15").Value = 24.359095
$ws.Range("N15").Value = 73.077285
$ws.Range("O15").Value = 0.9595205795530543
$ws.Range("P15").Value = 0.9595205795530543
$ws.Range("Q15").Value = 7.257150899081666
$ws.Range("R15").Value = 65.314358091735
$ws.Range("S15").Value = 0.0008999468001057623
$ws.Range("T15").Value = 0.0008999468001057625

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2979236666666666
$ws.Range("H16").Value = 0.893771
$ws.Range("I16").Value = 0.0009379129737112659
$ws.Range("J16").Value = 0.000937912973711266
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.6745613333333332
$ws.Range("N16").Value = 2.023684
$ws.Range("O16").Value = 0.02657140922124081
$ws.Range("P16").Value = 0.02657140922124081
$ws.Range("Q16").Value = 0.2009677858182222
$ws.Range("R16").Value = 1.808710072364
$ws.Range("S16").Value = 0.00002492166943839292
$ws.Range("T16").Value = 0.00002492166943839292
